$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 310
$ws1.Range("F3").Value = 1236
$ws1.Range("F4").Value = 17001
$ws1.Range("F6").Value = 1662
$ws1.Range("F9").Value = 13
$ws1.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202408/OZRx7O051723788701169.jpeg"
$ws1.Range("F11").Value = 230
$ws1.Range("F12").Value = 131
$ws1.Range("F13").Value = 11766
$ws1.Range("F15").Value = 22
$ws1.Range("F16").Value = 1448
$ws1.Range("F17").Value = 4676
$ws1.Range("F18").Value = 484
$ws1.Range("F19").Value = 19
$ws1.Range("F21").Value = 77
$ws1.Range("F22").Value = 909
$ws1.Range("F25").Value = 34
$ws1.Range("F26").Value = 5214

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 310
$ws4.Range("F4").Value = 1236
$ws4.Range("F5").Value = 17001
$ws4.Range("F7").Value = 1662
$ws4.Range("F10").Value = 13
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/OZRx7O051723788701169.jpeg"
$ws4.Range("F12").Value = 230
$ws4.Range("F13").Value = 131
$ws4.Range("F16").Value = 11766
$ws4.Range("F18").Value = 22
$ws4.Range("F19").Value = 1448
$ws4.Range("F20").Value = 4676
$ws4.Range("F21").Value = 484
$ws4.Range("F22").Value = 19
$ws4.Range("F24").Value = 77
$ws4.Range("F25").Value = 909
$ws4.Range("F28").Value = 34
$ws4.Range("F29").Value = 5214
